# ---------------------------------------------------------------------------
# Commit: "Wed, Jul 08, 2020  2:06:09 PM"
#
# Two logical changes are applied:
#
#  1. A table on slide 16 switches its table style (tableStyleId) from the
#     custom "{32237042-320B-40C1-947F-35EFD2679B8D}" style to the built-in
#     style "{041DBB69-4584-4AFA-A59A-4C8234FEE8A0}".
#
#  2. The deck's two theme parts (ppt/theme/theme1.xml = "Office Theme" and
#     ppt/theme/theme2.xml = "Integral") swap places: the theme actually
#     driving the slide master (today "Integral") becomes the default
#     "Office Theme" colour palette, and vice versa.  The font scheme and
#     format (fill/line/effect) scheme are already byte-identical between
#     the two theme parts, so only the 12 colour-scheme slots
#     (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink) actually need to change.
#     The only colour-scheme this automation surface can reach/write is the
#     one used by the slide master, so that is the one we update.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style id on the slide-16 table.
# ---------------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{041DBB69-4584-4AFA-A59A-4C8234FEE8A0}")
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the live theme's colour scheme from "Integral" to "Office Theme".
#    RGBColor.RGB uses the usual VBA RGB() encoding (R + G*256 + B*65536).
# ---------------------------------------------------------------------------
$officeThemeColors = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)

$masterColorScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $masterColorScheme.Colors($i).RGB = $officeThemeColors[$i - 1]
}
